# Auto update Excel log 2026-02-04 14:32:36
# Appends newly captured sensor readings to the PIR, Humidity and
# Temperature logs (all for the Bathroom sensor on 2026-02-04).

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        $ws,
        $startRow,
        $rows
    )

    $i = 0
    foreach ($r in $rows) {
        $rowNum = $startRow + $i

        # Column A (Date) - force text so "2026-02-04" is not
        # auto-converted into a date serial number.
        $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
        $ws.Cells.Item($rowNum, 1).Value = $r[0]
        $ws.Cells.Item($rowNum, 1).Style = "Normal"

        # Column B (Timestamp)
        $ws.Cells.Item($rowNum, 2).Value = $r[1]

        # Column C (Hour)
        $ws.Cells.Item($rowNum, 3).Value = $r[2]

        # Column D (Location)
        $ws.Cells.Item($rowNum, 4).Value = $r[3]

        # Column E (Value) - force text so percentage-looking values
        # such as "79.9%" are not auto-converted into numbers.
        $ws.Cells.Item($rowNum, 5).NumberFormat = "@"
        $ws.Cells.Item($rowNum, 5).Value = $r[4]
        $ws.Cells.Item($rowNum, 5).Style = "Normal"

        # Column F (Status)
        $ws.Cells.Item($rowNum, 6).Value = $r[5]

        $i = $i + 1
    }
}

# ---------------------------------------------------------------------
# PIR sheet - rows 366-378
# ---------------------------------------------------------------------
$pirRows = @(
    @("2026-02-04","14:31:33","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:31:35","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:31:40","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:31:45","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:31:50","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:31:55","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:32:00","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:32:09","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:32:10","14:00","Bathroom","Motion Detected","Active"),
    @("2026-02-04","14:32:16","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:32:21","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:32:26","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:32:27","14:00","Bathroom","Motion Detected","Active")
)
$wsPIR = $wb.Worksheets.Item("PIR")
Add-LogRows $wsPIR 366 $pirRows

# ---------------------------------------------------------------------
# Humidity sheet - rows 296-307
# ---------------------------------------------------------------------
$humidityRows = @(
    @("2026-02-04","14:31:31","14:00","Bathroom","79.9%","Active"),
    @("2026-02-04","14:31:34","14:00","Bathroom","79.1%","Active"),
    @("2026-02-04","14:31:37","14:00","Bathroom","80.0%","Active"),
    @("2026-02-04","14:31:42","14:00","Bathroom","79.1%","Active"),
    @("2026-02-04","14:31:47","14:00","Bathroom","80.0%","Active"),
    @("2026-02-04","14:31:52","14:00","Bathroom","79.1%","Active"),
    @("2026-02-04","14:31:57","14:00","Bathroom","79.9%","Active"),
    @("2026-02-04","14:32:02","14:00","Bathroom","78.9%","Active"),
    @("2026-02-04","14:32:07","14:00","Bathroom","80.0%","Active"),
    @("2026-02-04","14:32:12","14:00","Bathroom","78.9%","Active"),
    @("2026-02-04","14:32:17","14:00","Bathroom","79.9%","Active"),
    @("2026-02-04","14:32:23","14:00","Bathroom","78.9%","Active")
)
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRows $wsHumidity 296 $humidityRows

# ---------------------------------------------------------------------
# Temperature sheet - rows 296-307
# ---------------------------------------------------------------------
$temperatureRows = @(
    @("2026-02-04","14:31:32","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:31:35","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:31:38","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:31:43","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:31:48","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:31:53","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:31:58","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:32:03","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:32:08","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:32:13","14:00","Bathroom","24.2C","Active"),
    @("2026-02-04","14:32:18","14:00","Bathroom","24.3C","Active"),
    @("2026-02-04","14:32:23","14:00","Bathroom","24.3C","Active")
)
$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRows $wsTemperature 296 $temperatureRows
